# Sample.xlsx update: add an "automatic numbering" (TextName "^") example
# block as rows 27-30 on both worksheets, matching the upstream commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("English+Chinese")
$ws2 = $wb.Worksheets.Item("Japanese")

# --- make room: insert 4 new rows at row 27 on both sheets -----------------
$ws1.Range("A27:A30").EntireRow.Insert()
$ws2.Range("A27:A30").EntireRow.Insert()

# --- Sheet 1: "English+Chinese" --------------------------------------------
$ws1.Range("B27").Value = 'If the TextName is "^", it means that it inherits the prefix of the previous line TextName and automatically increments the number of the suffix. This is the automatic numbering function. This requires that the first of a series of automatic numbering must explicitly specify the prefix and initial suffix number. The advantage of using this function is that you don''t have to manually modify the TextName when inserting, deleting, or adjusting the order to ensure that the numbers are continuous. The program can also access these rows through a loop.'
$ws1.Range("C27").Value = 'TextName 如果是“^”则表示沿用前一行 TextName 的前缀并自动递增后缀的编号，这就是自动编号功能。这要求一系列自动编号的第一个必须明确指定前缀和初始后缀编号。使用这个功能的好处是，在对行进行插入、删除或者调整顺序时不必手动修改 TextName 来保证编号连续，程序也可以通过循环遍历的方式来访问这些行。'

# carry over the same formatting (wrap, font, style 10) used by the other
# "instructions" rows (e.g. row 23) onto the new row 27
$ws1.Range("B23:C23").Copy()
$ws1.Range("B27:C27").PasteSpecial(-4122)
$ws1.Rows.Item(27).RowHeight = 99.75

$ws1.Range("A28").Value = "Dialog0"
$ws1.Range("B28").Value = "This is the first in a series of automatic numbering."
$ws1.Range("C28").Value = "这是一系列自动编号的第一个。"

$ws1.Range("A29").Value = "^"
$ws1.Range("B29").Value = 'The final TextName of this line will be "Dialog1".'
$ws1.Range("C29").Value = "这一行的实际 TextName 会是“Dialog1”。"

$ws1.Range("A30").Value = "^"
$ws1.Range("B30").Value = 'The final TextName of this line will be "Dialog2".'
$ws1.Range("C30").Value = "这一行的实际 TextName 会是“Dialog2”。"

# --- Sheet 2: "Japanese" ----------------------------------------------------
$ws2.Range("B27").Value = 'If the TextName is "^", it means that it inherits the prefix of the previous line TextName and automatically increments the number of the suffix. This is the automatic numbering function. This requires that the first of a series of automatic numbering must explicitly specify the prefix and initial suffix number. The advantage of using this function is that you don''t have to manually modify the TextName when inserting, deleting, or adjusting the order to ensure that the numbers are continuous. The program can also access these rows through a loop.'
$ws2.Range("D27").Value = 'TextNameが「^」の場合、前の行のTextNameのプレフィックスを継承し、サフィックスの番号を自動的にインクリメントすることを意味します。これは自動番号付け機能です。 これには、一連の自動番号付けの最初で、プレフィックスと初期サフィックス番号を明示的に指定する必要があります。 この関数を使用する利点は、番号が連続するように順序を挿入、削除、または調整するときにTextNameを手動で変更する必要がないことです。また、プログラムはループを介してこれらの行にアクセスできます。'

$ws2.Range("B23").Copy()
$ws2.Range("B27").PasteSpecial(-4122)
$ws2.Range("D23").Copy()
$ws2.Range("D27").PasteSpecial(-4122)
$ws2.Rows.Item(27).RowHeight = 99.75

$ws2.Range("A28").Value = "Dialog0"
$ws2.Range("B28").Value = "This is the first in a series of automatic numbering."
$ws2.Range("D28").Value = "これは一連の自動番号付けの最初のものです。"

$ws2.Range("A29").Value = "^"
$ws2.Range("B29").Value = 'The final TextName of this line will be "Dialog1".'
$ws2.Range("D29").Value = "この行の最後のTextNameは「Dialog1」になります。"

$ws2.Range("A30").Value = "^"
$ws2.Range("B30").Value = 'The final TextName of this line will be "Dialog2".'
$ws2.Range("D30").Value = "この行の最後のTextNameは「Dialog2」になります。"

# --- cosmetic adjustments that shipped with the same commit ----------------
# (the host quantizes column width to 1/7-character steps, same as Excel's
# own pixel-snapping; 55.1 is the input that lands closest to the authored
# 55.875 stored width)
$ws1.Columns.Item(3).ColumnWidth = 55.1

# update the saved scroll position / selection on each sheet; re-select
# sheet 1 last so it stays the active ("tabSelected") tab, as in the source
$ws2.Application.ActiveWindow.ScrollRow = 4
$ws2.Range("D40").Select()

$ws1.Application.ActiveWindow.ScrollRow = 13
$ws1.Range("B28").Select()
